# Gallery_Sounder_FIM.xlsx edit:
#  - "Sounders" sheet renamed to "Germany"; its NGC code + Conventional
#    Sounders block updated to include the SB520 rows.
#  - New "Belgium" sheet added (copied from "Czech" so formatting/column
#    widths match) between "Germany" and "Czech", with its own market data.
#  - "Czech" sheet's NGC code + Conventional Sounders block updated the
#    same way as Germany/Belgium.
#  - "Czech" ends up the active/selected tab, as in the source workbook.
#
# NOTE: after Worksheet.Copy() inserts a sheet, previously-held sheet
# object variables can resolve to the wrong physical sheet (an engine
# quirk), so sheet references are re-fetched by name via
# Worksheets.Item(...) immediately after every Copy()/rename.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "Sounders" -> "Germany"
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Sounders").Name = "Germany"

# ---------------------------------------------------------------------
# 2. Create "Belgium" by copying "Czech" (keeps identical column widths /
#    styles) directly in front of "Czech", then rename the copy.
# ---------------------------------------------------------------------
$czechOrig = $wb.Worksheets.Item("Czech")
$czechOrig.Copy($czechOrig, $null)
$wb.Worksheets.Item("Czech (2)").Name = "Belgium"

# Re-fetch every sheet reference fresh now that the copy/rename is done.
$germany = $wb.Worksheets.Item("Germany")
$belgium = $wb.Worksheets.Item("Belgium")
$czech = $wb.Worksheets.Item("Czech")

# ---------------------------------------------------------------------
# 3. Update the "Germany" sheet content
# ---------------------------------------------------------------------
$germany.Range("B4").Value = "NGC-3475/T1763"

$germany.Range("A8").Value = "SB520"
$germany.Range("A9").Value = "SB520 c/w Front Cover"

# Rows 11/12 are new -- copy row 8's format down first, then fill values
$germany.Range("A8").Copy()
$germany.Range("A11:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$germany.Range("A10").Value = "Generic Sounder"
$germany.Range("A11").Value = "Wg"
$germany.Range("A12").Value = "Conventional Sounders"

$germany.Activate()
$germany.Range("A7").Select()

# ---------------------------------------------------------------------
# 4. Update the "Belgium" sheet content (new market data)
# ---------------------------------------------------------------------
$belgium.Range("B2").Value = "Belgium Market"
$belgium.Range("B4").Value = "NGC-3478/T2290/2326"

$belgium.Range("A8").Value = "SB520"
$belgium.Range("A9").Value = "SB520 c/w Front Cover"

# Rows 10/11/12 are new -- copy row 8's (now plain) format down, then fill
$belgium.Range("A8").Copy()
$belgium.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$belgium.Range("A10").Value = "Generic Sounder"
$belgium.Range("A11").Value = "Wg"
$belgium.Range("A12").Value = "Conventional Sounders"

$belgium.Activate()
$belgium.Range("B4").Select()

# ---------------------------------------------------------------------
# 5. Update the "Czech" sheet content
# ---------------------------------------------------------------------
$czech.Range("B4").Value = "NGC-3477/T1851/T1863"

$czech.Range("A8").Value = "SB520"
$czech.Range("A9").Value = "SB520 c/w Front Cover"

$czech.Range("A8").Copy()
$czech.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$czech.Range("A10").Value = "Generic Sounder"
$czech.Range("A11").Value = "Wg"
$czech.Range("A12").Value = "Conventional Sounders"

# Czech ends up the active sheet/tab, matching the source workbook.
$czech.Activate()
$czech.Range("B5").Select()
